# Mise à jour du classement - 31.03.2025 à 16:07

$wb = $excel.ActiveWorkbook

$updateText = "Dernière update le 31.03.25 à 16:07"

# --- leaderboard2 : "Qui a attrapé le plus de Cobblemons ?" ---
$ws1 = $wb.Worksheets.Item("leaderboard2")
$ws1.Range("D3").Value = 551
$ws1.Range("D4").Value = 513
$ws1.Range("B13").Value = $updateText

# --- leaderboard3 : "Qui a attrapé le plus de Shiny Cobblemons ?" ---
$ws2 = $wb.Worksheets.Item("leaderboard3")
$ws2.Range("D3").Value = 84
$ws2.Range("D4").Value = 44
$ws2.Range("B13").Value = $updateText

# --- leaderboard4 : "Qui a attrapé le plus de Cobblemons légendaires ?" ---
# Ranks 1 and 2 swap (ArtyumsM moves into 1st place, BKZRackham into 2nd),
# each with an incremented count.
$ws3 = $wb.Worksheets.Item("leaderboard4")
$ws3.Range("C3").Value = "ArtyumsM"
$ws3.Range("D3").Value = 9
$ws3.Range("C4").Value = "BKZRackham"
$ws3.Range("D4").Value = 8
$ws3.Range("B13").Value = $updateText
